$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 2-5: "Identificador" (D) switches from a running number to
# --- a text ID, and "Tipo" (E) switches from the number 1 to the text "Person".
$ws.Range("D2").Value = "59019237X"
$ws.Range("E2").Value = "Person"

$ws.Range("D3").Value = "57104958S"
$ws.Range("E3").Value = "Person"

$ws.Range("D4").Value = "10573947Y"
$ws.Range("E4").Value = "Person"

$ws.Range("D5").Value = "58429367Y"
$ws.Range("E5").Value = "Person"

# --- New row 6: Sensor1
$ws.Range("A6").Value = "Sensor1"
$ws.Range("B6").Value = "102.00.1"
$c6 = $ws.Range("C6")
$c6.Value = "sensor1@sensor.com"
$ws.Range("D6").Value = 3000
$ws.Range("E6").Value = "Sensor"

# --- New row 7: Sensor2
$ws.Range("A7").Value = "Sensor2"
$ws.Range("B7").Value = "999.1.120"
$c7 = $ws.Range("C7")
$c7.Value = "sensor2@sensor.com"
$ws.Range("D7").Value = 19800
$ws.Range("E7").Value = "Sensor"

# --- New row 8: Entity1
$ws.Range("A8").Value = "Entity1"
$ws.Range("B8").Value = "192.168.0.25"
$c8 = $ws.Range("C8")
$c8.Value = "entity@entity.com"
$ws.Range("D8").Value = "EntityX100"
$ws.Range("E8").Value = "Entity"

# --- Number formats: "Identificador" text format on the new sensor rows,
# --- then the header cell picks up a date-like number format (matches the
# --- source edit), in the same order the underlying styles were created.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D1").NumberFormat = "d-mmm-yy"

# --- Hyperlinks for the newly added contact e-mails.
$ws.Hyperlinks.Add($c6, "mailto:sensor1@sensor.com") | Out-Null
$c6.Style = "Hipervínculo"
$ws.Hyperlinks.Add($c7, "mailto:sensor2@sensor.com") | Out-Null
$c7.Style = "Hipervínculo"
$ws.Hyperlinks.Add($c8, "mailto:entity@entity.com") | Out-Null
$c8.Style = "Hipervínculo"

# --- Selection moves to D1, matching the saved cursor position.
$ws.Range("D1").Select() | Out-Null
